$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.19581585852348837,
    -0.0059999999758062472,
    -0.0039999999791273666,
    -0.007999999961564086,
    -0.0029999999789467324,
    -0.0019999999780857536,
    -0.0099999999463618039,
    -0.0099999999446653831,
    -0.0019999999737976282,
    -0.0019999999720337058,
    -0.0029999999680567768,
    -0.0034999999656872838,
    -0.0034999999640819013,
    -0.0025418039961495609,
    -0.00099999997295174836,
    -0.0019999999686599601,
    -0.0019999999681035163,
    -0.0039999999601914027,
    -0.035817198935969241,
    -0.0039999999819606558,
    -0.0039999999817750265,
    -0.0039999999816515697,
    -0.004999999973972713,
    -0.019999999913491884,
    -0.019999999912392319,
    -0.0024999999760595415,
    -0.0024999999756096791,
    -0.0019999999757107645,
    -0.0069999999550658387,
    -0.059999999749056343,
    -0.0069999999573884253,
    -0.0099999999467232925,
    0.061390397593225288
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
